$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: Insert 2 new rows at position 4 (shifts old rows 4-29 down to rows 6-31)
$ws.Rows("4:5").Insert()

# Step 2: Fix column A sequence numbers for shifted rows (6-31): they keep old value, need +2
for ($r = 6; $r -le 31; $r++) {
    $old = $ws.Cells.Item($r, 1).Value()
    $ws.Cells.Item($r, 1).Value = $old + 2
}

# Step 3: Rename "Thomas Hex" -> "Matthies Hex" wherever it occurs (now row 11, col B)
$ws.Cells.Item(11, 2).Value = "Matthies Hex"

# Step 4: Fill new row 4 (Holden) and row 5 (Rizzie Spiral)
$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "Holden"
$row4 = New-Object 'object[,]' 1,21
$row4vals = @(1.082600945368693, 1.082600945368693, 1.01350422820951, 1.012518769963506, 1.073287277478802, 0.9396312102174148, 0.8860049125673761, 0.959685783364458, 0.96012387523193, 0.8860049125673761, 1.082600945368693, 1.082600945368693, 0.96012387523193, 0.923064393899653, 0.9868140517207202, 0.976243244389333, 0.953211005336272, 0.976243244389333, 0.9855584903443773, 1.00496698134924, 0.9909196253002113)
for ($i = 0; $i -lt 21; $i++) { $row4[0,$i] = $row4vals[$i] }
$ws.Range("C4:W4").Value = $row4

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "Rizzie Spiral"
$row5 = New-Object 'object[,]' 1,21
$row5vals = @(1.85177977969318, 1.85177977969318, 0.8482465395264102, 1.032515094354373, 1.110791808156038, 0.9513023617946759, 1.27164900445721, 0.6615310458603386, 0.7816740777688083, 1.27164900445721, 1.85177977969318, 1.85177977969318, 0.7816740777688083, 1.026661541113009, 0.8149603086476092, 1.301700953973066, 0.967189873917476, 1.301700953973066, 1.188337350361402, 1.321025836227757, 1.063686213951379)
for ($i = 0; $i -lt 21; $i++) { $row5[0,$i] = $row5vals[$i] }
$ws.Range("C5:W5").Value = $row5

# Step 5: Fix styles for new rows 4,5 column A (match row 3 style - bold/border/centered)
$ws.Range("A4:A5").Font.Bold = $true
$ws.Range("A4:A5").HorizontalAlignment = -4108
$ws.Range("A4:A5").VerticalAlignment = -4160
$ws.Range("A4:A5").Borders.LineStyle = 1
